# Apply the commit's change to lifelink_db.xlsx:
#  1. Donors sheet: append a new donor record as row 16.
#  2. Patients sheet: drop the stray date-number-format style that was
#     applied to J10 (createdAt) so it matches the other rows (General).

$wb = $excel.ActiveWorkbook
$donors = $wb.Worksheets.Item("Donors")
$patients = $wb.Worksheets.Item("Patients")

# ---------------------------------------------------------------------
# 1. New Donors row (row 16)
# ---------------------------------------------------------------------
$donors.Range("A16").Value = "692a4edc07042cb8aa370cd8"
$donors.Range("B16").Value = "Ashwini Shenoy B"
$donors.Range("C16").Value = "ashenoyb@gmail.com"

# D16 (phone number) must stay text, like the other phone-number cells
# in this column, not be auto-coerced into a number. Enter it with a
# leading quote so Excel keeps it as text, then re-apply the plain
# General formatting from a neighboring cell so the cell itself doesn't
# keep a one-off "quote prefix" style.
$donors.Range("D16").Value = "'7026438371"
$donors.Range("D15").Copy()
$donors.Range("D16").PasteSpecial(-4122)

$donors.Range("E16").Value = "O+"
$donors.Range("F16").Value = "My Current Location"
$donors.Range("G16").Value = 75.1239547
$donors.Range("H16").Value = 15.3647083
$donors.Range("I16").Value = $false
$donors.Range("J16").Value = $false
$donors.Range("K16").Value = $true
$donors.Range("L16").Value = 45990.29850224537

# L16 (createdAt) should carry the same date number format used by the
# other "createdAt" style already present in this workbook (style index
# 1, numFmtId 14). Copy that formatting over instead of typing a format
# string so we reuse the existing style rather than registering a new one.
# NOTE: this must run before Patients!J10's own format is reset below,
# since Patients!J10 is the only other cell currently holding that style.
$patients.Range("J10").Copy()
$donors.Range("L16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Patients!J10 loses its one-off date format (back to General, like
#    every other createdAt cell in that column).
# ---------------------------------------------------------------------
$patients.Range("J9").Copy()
$patients.Range("J10").PasteSpecial(-4122)
# Re-assert the value in case the format paste touched it.
$patients.Range("J10").Value = 45989.91622431713

Write-Host "Applied lifelink_db.xlsx donor row + style update"
